{"js": "// CV update: move masthead date forward, add a new \"Amazon Scholar\" bullet,\n// and tighten the existing \"Principal Economist\" bullet to reflect that the\n// Amazon role ended in 2024 and the group is now called \"SEAS\".\n\nconst body = context.document.body;\n\n// 1) Masthead date: \"July 2024\" -> \"September 2024\"\nconst julySearch = body.search(\"July\", { matchCase: true, matchWholeWord: true });\njulySearch.load(\"items\");\nawait context.sync();\nif (julySearch.items.length > 0) {\n  julySearch.items[0].insertText(\"September\", \"Replace\");\n  await context.sync();\n}\n\n// 2) Find the existing \"Principal Economist, Amazon.com, ...\" bullet\n//    paragraph and insert a brand-new bullet paragraph right before it for\n//    the new \"Amazon Scholar\" role (numbering / font formatting is\n//    inherited automatically from the paragraph it is inserted next to).\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet targetPara = null;\nfor (const p of paragraphs.items) {\n  if (p.text.indexOf(\"Principal Economist, Amazon.com\") === 0) {\n    targetPara = p;\n    break;\n  }\n}\n\ntargetPara.insertParagraph(\n  \"Amazon Scholar, Amazon.com, Stores Economics and Science Group (SEAS), 2024-Present\",\n  \"Before\"\n);\nawait context.sync();\n\n// 3) Update the original \"Principal Economist\" bullet:\n//    - \"(formerly Core AI)\" -> \"(SEAS, formerly Core AI)\"\n//    - trailing \"2022-Present\" -> \"2022-2024\"\nconst formerlySearch = body.search(\"Stores Economics and Science Group (formerly\", { matchCase: true });\nformerlySearch.load(\"items\");\nawait context.sync();\nif (formerlySearch.items.length > 0) {\n  formerlySearch.items[0].insertText(\n    \"Stores Economics and Science Group (SEAS, formerly\",\n    \"Replace\"\n  );\n  await context.sync();\n}\n\nconst presentSearch = body.search(\"Core AI), 2022-Present\", { matchCase: true });\npresentSearch.load(\"items\");\nawait context.sync();\nif (presentSearch.items.length > 0) {\n  presentSearch.items[0].insertText(\"Core AI), 2022-2024\", \"Replace\");\n  await context.sync();\n}\n", "ps1": "# CV update: move masthead date forward, add a new \"Amazon Scholar\" bullet,\n# and tighten the existing \"Principal Economist\" bullet to reflect that the\n# Amazon role ended in 2024 and the group is now called \"SEAS\".\n\n$doc = $word.ActiveDocument\n\n# 1) Masthead date: \"July 2024\" -> \"September 2024\"\n$find = $doc.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"July\", $false, $false, $false, $false, $false, $true, 1, $false, \"September\", 2)\n\n# 2) Find the existing \"Principal Economist, Amazon.com, ...\" bullet paragraph\n#    and insert a brand-new bullet paragraph right before it for the new\n#    \"Amazon Scholar\" role (same list / formatting is inherited automatically).\n$targetPara = $null\nforeach ($p in $doc.Paragraphs) {\n    if ($p.Range.Text -like \"Principal Economist, Amazon.com*\") {\n        $targetPara = $p\n        break\n    }\n}\n\n$targetPara.Range.InsertParagraphBefore()\n\n# InsertParagraphBefore() splits off a new, empty paragraph in front of the\n# original text and re-seats $targetPara onto that new (still empty)\n# paragraph, so we can set its text directly.\n$targetPara.Range.Text = \"Amazon Scholar, Amazon.com, Stores Economics and Science Group (SEAS), 2024-Present\"\n\n# 3) Update the original \"Principal Economist\" bullet:\n#    - \"(formerly Core AI)\" -> \"(SEAS, formerly Core AI)\"\n#    - trailing \"2022-Present\" -> \"2022-2024\"\n$find2 = $doc.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Execute(\"Stores Economics and Science Group (formerly\", $false, $false, $false, $false, $false, $true, 1, $false, \"Stores Economics and Science Group (SEAS, formerly\", 2)\n\n$find3 = $doc.Content.Find\n$find3.ClearFormatting()\n$find3.Replacement.ClearFormatting()\n$find3.Execute(\"Core AI), 2022-Present\", $false, $false, $false, $false, $false, $true, 1, $false, \"Core AI), 2022-2024\", 2)\n"}
